$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values after repulling data / recalculating means.
$updates = @{
    2  = 1
    3  = -2
    7  = 3
    10 = -1
    11 = -2
    18 = 1
    21 = 0
    28 = -1
    33 = -3
    36 = -2
    37 = -1
    41 = -2
    42 = -2
    48 = 2
    50 = -2
    56 = -1
    57 = -3
    59 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
